# MOSIP Features Roadmap update
# Adds new tracked-requirement rows (151-155 / sheet rows 154-158) to the
# MOSIP_Feature_Roadmap sheet, covering new IDA / ID-Repo / VID items that
# were reviewed & approved by Sasi/Ramesh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_Feature_Roadmap")

function Set-RoadmapRow {
    param(
        [int]$Row,
        [int]$SNo,
        [string]$Jira,
        [int]$Year,
        [int]$Month,
        [int]$Day,
        [string]$Reference,
        [string]$Module,
        [string]$ChangeType,
        [string]$ChangeDescription,
        [string]$Approver
    )

    $date = (Get-Date -Year $Year -Month $Month -Day $Day).Date
    $addrA = "A" + $Row
    $addrB = "B" + $Row
    $addrC = "C" + $Row
    $addrD = "D" + $Row
    $addrE = "E" + $Row
    $addrF = "F" + $Row
    $addrG = "G" + $Row
    $addrL = "L" + $Row
    $addrM = "M" + $Row
    $addrN = "N" + $Row
    $addrO = "O" + $Row

    $ws.Range($addrA).Value = $SNo
    $ws.Range($addrB).Value = $Jira
    $ws.Range($addrC).Value = $date
    $ws.Range($addrC).NumberFormat = "d-mmm-yy"
    $ws.Range($addrD).Value = $Reference
    $ws.Range($addrE).Value = $Module
    $ws.Range($addrF).Value = $ChangeType
    $ws.Range($addrG).Value = $ChangeDescription
    $ws.Range($addrL).Value = 1
    $ws.Range($addrM).Value = "Approved"
    $ws.Range($addrN).Value = $Approver
    $ws.Range($addrO).Value = $date
    $ws.Range($addrO).NumberFormat = "d-mmm-yy"
}

Set-RoadmapRow 154 151 "MOS-21582" 2019 3 26 `
    "API Specification Changes for IDA based on MDS review by Sasi/Ramesh" `
    "ID-Authentication" "New" `
    "Additional or Modification of attributes in API Specs based on review " `
    "Ramesh"

Set-RoadmapRow 155 152 "MOS-21583" 2019 3 26 `
    "Design Change of ID-Repo based on Security review by Sasi/Ramesh" `
    "ID-Authentication" "New" `
    "Design Change of ID-Repo based on Security review by Sasi/Ramesh" `
    "Ramesh"

Set-RoadmapRow 156 153 "MOS-21584" 2019 3 26 `
    "Design Change of IDA based on Security review by Sasi/Ramesh" `
    "ID-Authentication" "New" `
    "Design Change of IDA based on Security review by Sasi/Ramesh" `
    "Ramesh"

Set-RoadmapRow 157 154 "MOS-21585" 2019 4 1 `
    "Mapping of platform address attributes in IDA based on Morrocco Address Structure" `
    "ID-Authentication" "New" `
    "Mapping of platform address attributes in IDA based on Morrocco Address Structure" `
    "Shrikant"

Set-RoadmapRow 158 155 "MOS-21327" 2019 3 28 `
    "Integrate with new VID Generator API" `
    "ID-Authentication" "New" `
    "Integrate with the new VID generator component based on the VID policy/type defined" `
    "Ramesh"
